# Update win-probability matrix cells with refreshed simulation results
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 0.21
$ws.Range("C2").Value = 0.53
$ws.Range("J2").Value = 0.02
$ws.Range("P2").Value = 0.13
$ws.Range("S2").Value = 0.11

# Row 3
$ws.Range("C3").Value = 0.03636363636363636
$ws.Range("J3").Value = 0.01818181818181818
$ws.Range("P3").Value = 0.6909090909090909
$ws.Range("S3").Value = 0.2545454545454545

# Row 4
$ws.Range("J4").Value = 0.06666666666666667
$ws.Range("P4").Value = 0.7333333333333333
$ws.Range("S4").Value = 0.2

# Row 6
$ws.Range("B6").Value = 0.06315789473684211
$ws.Range("D6").Value = 0.01052631578947368
$ws.Range("F6").Value = 0.09473684210526316
$ws.Range("J6").Value = 0.1052631578947368
$ws.Range("O6").Value = 0.02105263157894737
$ws.Range("Q6").Value = 0.1578947368421053
$ws.Range("R6").Value = 0.07368421052631578
$ws.Range("S6").Value = 0.4736842105263158

# Row 7
$ws.Range("B7").Value = 0.08955223880597014
$ws.Range("D7").Value = 0.02985074626865672
$ws.Range("F7").Value = 0.04477611940298507
$ws.Range("J7").Value = 0.08955223880597014
$ws.Range("O7").Value = 0.04477611940298507
$ws.Range("Q7").Value = 0.2835820895522388
$ws.Range("R7").Value = 0.08955223880597014
$ws.Range("S7").Value = 0.3283582089552239

# Row 8
$ws.Range("B8").Value = 0.06
$ws.Range("D8").Value = 0.01
$ws.Range("F8").Value = 0.115
$ws.Range("J8").Value = 0.15
$ws.Range("O8").Value = 0.015
$ws.Range("Q8").Value = 0.19
$ws.Range("R8").Value = 0.105
$ws.Range("S8").Value = 0.355

# Row 9
$ws.Range("B9").Value = 0.09302325581395349
$ws.Range("D9").Value = 0.02325581395348837
$ws.Range("F9").Value = 0.04651162790697674
$ws.Range("J9").Value = 0.06976744186046512
$ws.Range("O9").Value = 0.02325581395348837
$ws.Range("Q9").Value = 0.2325581395348837
$ws.Range("R9").Value = 0.04651162790697674
$ws.Range("S9").Value = 0.4651162790697674

# Row 10
$ws.Range("B10").Value = 0.08918406072106262
$ws.Range("D10").Value = 0.01518026565464896
$ws.Range("F10").Value = 0.06831119544592031
$ws.Range("J10").Value = 0.127134724857685
$ws.Range("O10").Value = 0.02466793168880456
$ws.Range("Q10").Value = 0.1878557874762808
$ws.Range("R10").Value = 0.1024667931688805
$ws.Range("S10").Value = 0.3851992409867173

# Row 11
$ws.Range("G11").Value = 0.1157894736842105
$ws.Range("J11").Value = 0.1157894736842105
$ws.Range("K11").Value = 0.1578947368421053
$ws.Range("L11").Value = 0.6
$ws.Range("S11").Value = 0.01052631578947368

# Row 12
$ws.Range("G12").Value = 0.7457627118644068
$ws.Range("J12").Value = 0.2033898305084746
$ws.Range("L12").Value = 0.01694915254237288
$ws.Range("S12").Value = 0.03389830508474576

# Row 13
$ws.Range("G13").Value = 0.8235294117647058
$ws.Range("J13").Value = 0.1764705882352941

# Row 15
$ws.Range("F15").Value = 0.009433962264150943
$ws.Range("H15").Value = 0.160377358490566
$ws.Range("I15").Value = 0.0660377358490566
$ws.Range("J15").Value = 0.3962264150943396
$ws.Range("K15").Value = 0.05660377358490566
$ws.Range("M15").Value = 0.009433962264150943
$ws.Range("O15").Value = 0.1037735849056604
$ws.Range("S15").Value = 0.1981132075471698

# Row 16
$ws.Range("F16").Value = 0.01724137931034483
$ws.Range("H16").Value = 0.1896551724137931
$ws.Range("I16").Value = 0.05172413793103448
$ws.Range("J16").Value = 0.3448275862068966
$ws.Range("K16").Value = 0.1896551724137931
$ws.Range("M16").Value = 0.01724137931034483
$ws.Range("O16").Value = 0.08620689655172414
$ws.Range("S16").Value = 0.103448275862069

# Row 17
$ws.Range("F17").Value = 0.005208333333333333
$ws.Range("H17").Value = 0.21875
$ws.Range("I17").Value = 0.1041666666666667
$ws.Range("J17").Value = 0.4114583333333333
$ws.Range("K17").Value = 0.078125
$ws.Range("M17").Value = 0.015625
$ws.Range("O17").Value = 0.0625
$ws.Range("S17").Value = 0.1041666666666667

# Row 18
$ws.Range("H18").Value = 0.2043010752688172
$ws.Range("I18").Value = 0.08602150537634409
$ws.Range("J18").Value = 0.4946236559139785
$ws.Range("K18").Value = 0.05376344086021505
$ws.Range("M18").Value = 0.02150537634408602
$ws.Range("O18").Value = 0.07526881720430108
$ws.Range("S18").Value = 0.06451612903225806

# Row 19
$ws.Range("F19").Value = 0.0119047619047619
$ws.Range("H19").Value = 0.2261904761904762
$ws.Range("I19").Value = 0.09523809523809523
$ws.Range("J19").Value = 0.3928571428571428
$ws.Range("K19").Value = 0.08333333333333333
$ws.Range("M19").Value = 0.01984126984126984
$ws.Range("O19").Value = 0.07738095238095238
$ws.Range("S19").Value = 0.09325396825396826

Write-Host "Applied matrix updates"
